$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump version and update date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# --- Elements sheet: the ele-1/ext-1 invariant belongs to the Extension
#     element itself (row 1), not just Extension.extension (row 3) ---
$invariant = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ1").Value = $invariant
